# Fix font in some paragraph styles: "Bibliography" and "Footnote Text"
# previously had empty run properties (w:rPr/>); give them explicit
# Times New Roman fonts (and, for Bibliography, a 14pt size) to match
# the rest of the document's Times New Roman styling.

$d = $word.ActiveDocument

$bibliography = $d.Styles("Bibliography")
$bibliography.Font.Name = "Times New Roman"
$bibliography.Font.Size = 14

$footnoteText = $d.Styles("Footnote Text")
$footnoteText.Font.Name = "Times New Roman"
